$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.57620000000002
$ws.Range("A21").Value = -21.17510000000001
$ws.Range("A23").Value = -21.42440000000003
$ws.Range("A25").Value = -22.41280000000003
